$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows appended to the feed logs sheet (rows 96 and 97)
$ws.Range("A96").Value = 95
$ws.Range("B96").Value = 1
$ws.Range("C96").Value = "2024-06-16 22:13:29"
$ws.Range("D96").Value = 200
$ws.Range("E96").Value = 7

$ws.Range("A97").Value = 96
$ws.Range("B97").Value = 2
$ws.Range("C97").Value = "2024-06-16 22:13:29"
$ws.Range("D97").Value = 200
$ws.Range("E97").Value = 0
